# Rewrite the cash-flow table (rows 2-20) to match the updated financial + enrollment data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bordered/bold/centered header style used by column A (school name) down to the new rows.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A13:A20").PasteSpecial(-4122) | Out-Null

# Row 2: ARIZONA_STATE_UNIVERSITY
$row2 = New-Object "object[,]" 1,16
$row2[0,0] = "ARIZONA_STATE_UNIVERSITY"
$row2[0,1] = "2023‑2024"
$row2[0,2] = 194089000
$row2[0,3] = 220984000
$row2[0,4] = ""
$row2[0,5] = ""
$row2[0,6] = -655833
$row2[0,7] = 466829
$row2[0,8] = ""
$row2[0,9] = 54813
$row2[0,10] = 342354
$row2[0,11] = 120964
$row2[0,12] = ""
$row2[0,13] = ""
$row2[0,14] = 679424
$row2[0,15] = 78400
$ws.Range("A2:P2").Value = $row2

# Row 3: BRADLEY_UNIVERSITY
$row3 = New-Object "object[,]" 1,16
$row3[0,0] = "BRADLEY_UNIVERSITY"
$row3[0,1] = "2023‑2024"
$row3[0,2] = 17930
$row3[0,3] = 13560
$row3[0,4] = -4826
$row3[0,5] = ""
$row3[0,6] = -28229
$row3[0,7] = 2475
$row3[0,8] = ""
$row3[0,9] = 31106
$row3[0,10] = 16990
$row3[0,11] = 17335
$row3[0,12] = -4345
$row3[0,13] = ""
$row3[0,14] = 2865
$row3[0,15] = 5742
$ws.Range("A3:P3").Value = $row3

# Row 4: CALIFORNIA_STATE_UNIVERSITY
$row4 = New-Object "object[,]" 1,16
$row4[0,0] = "CALIFORNIA_STATE_UNIVERSITY"
$row4[0,1] = "2023‑2024"
$row4[0,2] = -502289
$row4[0,3] = 748772
$row4[0,4] = ""
$row4[0,5] = ""
$row4[0,6] = -7022740
$row4[0,7] = 1113761
$row4[0,8] = ""
$row4[0,9] = 1149921
$row4[0,10] = 47245
$row4[0,11] = 365998
$row4[0,12] = 640013
$row4[0,13] = ""
$row4[0,14] = 6764422
$row4[0,15] = 38789
$ws.Range("A4:P4").Value = $row4

# Row 5: CORNELL_UNIVERSITY
$row5 = New-Object "object[,]" 1,16
$row5[0,0] = "CORNELL_UNIVERSITY"
$row5[0,1] = "2023‑2024"
$row5[0,2] = 488140
$row5[0,3] = 293011
$row5[0,4] = ""
$row5[0,5] = ""
$row5[0,6] = -238803
$row5[0,7] = 504143
$row5[0,8] = ""
$row5[0,9] = -464863
$row5[0,10] = 1269517
$row5[0,11] = 738394
$row5[0,12] = ""
$row5[0,13] = ""
$row5[0,14] = 894914
$row5[0,15] = 191248
$ws.Range("A5:P5").Value = $row5

# Row 6: CULINARY_INSTITUTE_OF_AMERICA_T
$row6 = New-Object "object[,]" 1,16
$row6[0,0] = "CULINARY_INSTITUTE_OF_AMERICA_T"
$row6[0,1] = "2023‑2024"
$row6[0,2] = 39322692
$row6[0,3] = 10112852
$row6[0,4] = ""
$row6[0,5] = ""
$row6[0,6] = 19904473
$row6[0,7] = 11579120
$row6[0,8] = ""
$row6[0,9] = -27266765
$row6[0,10] = ""
$row6[0,11] = 5289884
$row6[0,12] = ""
$row6[0,13] = ""
$row6[0,14] = 328207
$row6[0,15] = -7034085
$ws.Range("A6:P6").Value = $row6

# Row 7: GANNON_UNIVERSITY
$row7 = New-Object "object[,]" 1,16
$row7[0,0] = "GANNON_UNIVERSITY"
$row7[0,1] = "2023‑2024"
$row7[0,2] = 8395
$row7[0,3] = 7401
$row7[0,4] = 6555
$row7[0,5] = ""
$row7[0,6] = 12813
$row7[0,7] = 10112
$row7[0,8] = ""
$row7[0,9] = -9060
$row7[0,10] = ""
$row7[0,11] = 2131
$row7[0,12] = -2402
$row7[0,13] = 1339
$row7[0,14] = -792
$row7[0,15] = 2961
$ws.Range("A7:P7").Value = $row7

# Row 8: LEWIS_UNIVERSITY
$row8 = New-Object "object[,]" 1,16
$row8[0,0] = "LEWIS_UNIVERSITY"
$row8[0,1] = "2023‑2024"
$row8[0,2] = 23471502
$row8[0,3] = 12280.064
$row8[0,4] = -4208019
$row8[0,5] = ""
$row8[0,6] = 7320786
$row8[0,7] = 18103845
$row8[0,8] = ""
$row8[0,9] = -13716182
$row8[0,10] = 550000
$row8[0,11] = 2599642
$row8[0,12] = ""
$row8[0,13] = 246841
$row8[0,14] = -1475420
$row8[0,15] = -7870816
$ws.Range("A8:P8").Value = $row8

# Row 9: MICHIGAN_STATE_UNIVERSITY
$row9 = New-Object "object[,]" 1,16
$row9[0,0] = "MICHIGAN_STATE_UNIVERSITY"
$row9[0,1] = "2023‑2024"
$row9[0,2] = 342054
$row9[0,3] = 238893
$row9[0,4] = ""
$row9[0,5] = ""
$row9[0,6] = -537531
$row9[0,7] = 328895
$row9[0,8] = ""
$row9[0,9] = 207987
$row9[0,10] = 466397
$row9[0,11] = 190389
$row9[0,12] = 367558
$row9[0,13] = ""
$row9[0,14] = 547395
$row9[0,15] = 766454
$ws.Range("A9:P9").Value = $row9

# Row 10: MOLLOY_COLLEGE
$row10 = New-Object "object[,]" 1,16
$row10[0,0] = "MOLLOY_COLLEGE"
$row10[0,1] = "2023‑2024"
$row10[0,2] = 6573828
$row10[0,3] = 8078552
$row10[0,4] = ""
$row10[0,5] = ""
$row10[0,6] = 14016100
$row10[0,7] = 2932101
$row10[0,8] = ""
$row10[0,9] = 15262911
$row10[0,10] = ""
$row10[0,11] = 2410000
$row10[0,12] = ""
$row10[0,13] = ""
$row10[0,14] = -2120213
$row10[0,15] = 27158798
$ws.Range("A10:P10").Value = $row10

# Row 11: MOUNT_ST_MARY_S_UNIVERSITY_INC
$row11 = New-Object "object[,]" 1,16
$row11[0,0] = "MOUNT_ST_MARY_S_UNIVERSITY_INC"
$row11[0,1] = "2023‑2024"
$row11[0,2] = 6951753
$row11[0,3] = 7514013
$row11[0,4] = ""
$row11[0,5] = ""
$row11[0,6] = 6297553
$row11[0,7] = 17662886
$row11[0,8] = ""
$row11[0,9] = -14596047
$row11[0,10] = ""
$row11[0,11] = 1325000
$row11[0,12] = ""
$row11[0,13] = ""
$row11[0,14] = 7533886
$row11[0,15] = -764608
$ws.Range("A11:P11").Value = $row11

# Row 12: NEW_YORK_UNIVERSITY
$row12 = New-Object "object[,]" 1,16
$row12[0,0] = "NEW_YORK_UNIVERSITY"
$row12[0,1] = "2023‑2024"
$row12[0,2] = 1608069
$row12[0,3] = 1128846
$row12[0,4] = ""
$row12[0,5] = ""
$row12[0,6] = 1333220
$row12[0,7] = 1724357
$row12[0,8] = ""
$row12[0,9] = -202895
$row12[0,10] = ""
$row12[0,11] = 225521
$row12[0,12] = ""
$row12[0,13] = ""
$row12[0,14] = 196971
$row12[0,15] = 1327296
$ws.Range("A12:P12").Value = $row12

# Row 13: OHIO_STATE_UNIVERSITY_THE
$row13 = New-Object "object[,]" 1,16
$row13[0,0] = "OHIO_STATE_UNIVERSITY_THE"
$row13[0,1] = "2023‑2024"
$row13[0,2] = 873349
$row13[0,3] = 616748
$row13[0,4] = ""
$row13[0,5] = ""
$row13[0,6] = ""
$row13[0,7] = ""
$row13[0,8] = ""
$row13[0,9] = ""
$row13[0,10] = ""
$row13[0,11] = ""
$row13[0,12] = ""
$row13[0,13] = ""
$row13[0,14] = ""
$row13[0,15] = ""
$ws.Range("A13:P13").Value = $row13

# Row 14: PRESIDENT___FELLOWS_OF_HARVARD_
$row14 = New-Object "object[,]" 1,16
$row14[0,0] = "PRESIDENT___FELLOWS_OF_HARVARD_"
$row14[0,1] = "2023‑2024"
$row14[0,2] = 2811779
$row14[0,3] = 440257
$row14[0,4] = ""
$row14[0,5] = ""
$row14[0,6] = -2151998
$row14[0,7] = 777016
$row14[0,8] = ""
$row14[0,9] = 1152481
$row14[0,10] = 2287902
$row14[0,11] = 1295556.2
$row14[0,12] = ""
$row14[0,13] = ""
$row14[0,14] = 1371336
$row14[0,15] = 371819
$ws.Range("A14:P14").Value = $row14

# Row 15: STEVENSON_UNIVERSITY_INC
$row15 = New-Object "object[,]" 1,16
$row15[0,0] = "STEVENSON_UNIVERSITY_INC"
$row15[0,1] = "2023‑2024"
$row15[0,2] = 440212
$row15[0,3] = 8079347
$row15[0,4] = 4220760
$row15[0,5] = 341011
$row15[0,6] = 2761885
$row15[0,7] = 10677588
$row15[0,8] = -8900915
$row15[0,9] = -11877590
$row15[0,10] = 10163015
$row15[0,11] = 2161567
$row15[0,12] = ""
$row15[0,13] = 1246407
$row15[0,14] = 9147855
$row15[0,15] = 32150
$ws.Range("A15:P15").Value = $row15

# Row 16: STEVENS_INSTITUTE_OF_TECHNOLOGY
$row16 = New-Object "object[,]" 1,16
$row16[0,0] = "STEVENS_INSTITUTE_OF_TECHNOLOGY"
$row16[0,1] = "2023‑2024"
$row16[0,2] = 57046
$row16[0,3] = ""
$row16[0,4] = ""
$row16[0,5] = ""
$row16[0,6] = 30949
$row16[0,7] = 44103
$row16[0,8] = ""
$row16[0,9] = -54402
$row16[0,10] = 828
$row16[0,11] = 6065
$row16[0,12] = ""
$row16[0,13] = ""
$row16[0,14] = 6218
$row16[0,15] = -17235
$ws.Range("A16:P16").Value = $row16

# Row 17: ST_LOUIS_UNIVERSITY_US
$row17 = New-Object "object[,]" 1,16
$row17[0,0] = "ST_LOUIS_UNIVERSITY_US"
$row17[0,1] = "2023‑2024"
$row17[0,2] = 48539
$row17[0,3] = 45517
$row17[0,4] = ""
$row17[0,5] = -11073
$row17[0,6] = -28308
$row17[0,7] = 48839
$row17[0,8] = ""
$row17[0,9] = -18771
$row17[0,10] = ""
$row17[0,11] = 11431
$row17[0,12] = ""
$row17[0,13] = ""
$row17[0,14] = 5335
$row17[0,15] = -41744
$ws.Range("A17:P17").Value = $row17

# Row 18: TEXAS_A_M_UNIVERSITY
$row18 = New-Object "object[,]" 1,16
$row18[0,0] = "TEXAS_A_M_UNIVERSITY"
$row18[0,1] = "2023‑2024"
$row18[0,2] = 1206282105.94
$row18[0,3] = 586407632.0700001
$row18[0,4] = ""
$row18[0,5] = ""
$row18[0,6] = -3091426326.9
$row18[0,7] = 920150736.05
$row18[0,8] = ""
$row18[0,9] = 277286556.04
$row18[0,10] = 956820491.0599999
$row18[0,11] = 613865951.75
$row18[0,12] = 157410000
$row18[0,13] = ""
$row18[0,14] = 3277592772.9
$row18[0,15] = 463453002.04
$ws.Range("A18:P18").Value = $row18

# Row 19: UNIVERSITY_OF_COLORADO
$row19 = New-Object "object[,]" 1,16
$row19[0,0] = "UNIVERSITY_OF_COLORADO"
$row19[0,1] = "2023‑2024"
$row19[0,2] = ""
$row19[0,3] = ""
$row19[0,4] = ""
$row19[0,5] = ""
$row19[0,6] = ""
$row19[0,7] = ""
$row19[0,8] = ""
$row19[0,9] = ""
$row19[0,10] = ""
$row19[0,11] = ""
$row19[0,12] = ""
$row19[0,13] = ""
$row19[0,14] = ""
$row19[0,15] = ""
$ws.Range("A19:P19").Value = $row19

# Row 20: UNIVERSITY_OF_MINNESOTA
$row20 = New-Object "object[,]" 1,16
$row20[0,0] = "UNIVERSITY_OF_MINNESOTA"
$row20[0,1] = "2023‑2024"
$row20[0,2] = 215488
$row20[0,3] = 258080
$row20[0,4] = ""
$row20[0,5] = ""
$row20[0,6] = -1345762
$row20[0,7] = 216236
$row20[0,8] = ""
$row20[0,9] = 167725
$row20[0,10] = 250552
$row20[0,11] = 627260
$row20[0,12] = -21785
$row20[0,13] = ""
$row20[0,14] = ""
$row20[0,15] = -60232
$ws.Range("A20:P20").Value = $row20

